$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99 (pushes existing rows 99.. down by one,
# and extends the sheet dimension from T147 to T148 automatically).
$ws.Rows(99).Insert()

# Populate the newly inserted row 99 with the new record.
$ws.Range("A99").Value = 3
$ws.Range("B99").Value = "Femacal de La Calera"
$ws.Range("C99").Value = "Coquimbo"
$ws.Range("D99").Value = 44510
$ws.Range("E99").Value = 5
$ws.Range("F99").Value = "Fruta"
$ws.Range("G99").Value = 100101
$ws.Range("H99").Value = "Berries"
$ws.Range("I99").Value = 100101001
$ws.Range("J99").Value = "Arándano (blue)"
$ws.Range("K99").Value = "Sin especificar"
$ws.Range("L99").Value = "Primera"
$ws.Range("M99").Value = 50
$ws.Range("N99").Value = 10000
$ws.Range("O99").Value = 10000
$ws.Range("P99").Value = 10000
$ws.Range("Q99").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R99").Value = "Provincia de Quillota"
$ws.Range("S99").Value = 6667
$ws.Range("T99").Value = 1.5
